$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1:J1").FormulaArray = '=CELL("width")'
$ws.Range("I2").Select()
